$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4314.778
$ws.Range("I32").Value = 3925.5
$ws.Range("J32").Value = 4626.2
$ws.Range("K32").Value = 3925.5
$ws.Range("L32").Value = 4626.2
$ws.Range("M32").Value = -3599.5
$ws.Range("N32").Value = -5278.2

$ws.Range("H116").Value = 224719.16
$ws.Range("I116").Value = 104040.836
$ws.Range("J116").Value = 280416.84
$ws.Range("K116").Value = 104040.836
$ws.Range("L116").Value = 280416.84
$ws.Range("M116").Value = -100598.836
$ws.Range("N116").Value = -287300.84

$ws.Range("H132").Value = 61789.848
$ws.Range("I132").Value = 68157.31
$ws.Range("J132").Value = 6074.5
$ws.Range("K132").Value = 204471.93
$ws.Range("L132").Value = 18223.5
$ws.Range("M132").Value = -201941.93
$ws.Range("N132").Value = -23283.5

$ws.Range("H137").Value = 1424.09
$ws.Range("J137").Value = 1259.0548
$ws.Range("L137").Value = 3777.1644
$ws.Range("N137").Value = -8877.1644

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 2028.8572
$ws.Range("I33").Value = 2028.8572
$ws.Range("K33").Value = 2028.8572
$ws.Range("M33").Value = -1699.8572

$ws.Range("H37").Value = 67252.87
$ws.Range("I37").Value = 38748.5
$ws.Range("J37").Value = 77618.09
$ws.Range("K37").Value = 38748.5
$ws.Range("L37").Value = 77618.09
$ws.Range("M37").Value = -38475.5
$ws.Range("N37").Value = -78164.09

$ws.Range("H45").Value = 2191.5334
$ws.Range("I45").Value = 2190.0833
$ws.Range("K45").Value = 2190.0833
$ws.Range("M45").Value = -1813.0833

$ws.Range("H74").Value = 11247.3125
$ws.Range("I74").Value = 13495.7
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 13495.7
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = -12621.7
$ws.Range("N74").Value = -9248

$ws.Range("H77").Value = 11247.3125
$ws.Range("I77").Value = 13495.7
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 67478.5
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -63110.5
$ws.Range("N77").Value = -46236

$ws.Range("H97").Value = 1803.5555
$ws.Range("I97").Value = 1702
$ws.Range("K97").Value = 1702
$ws.Range("M97").Value = -1206

$ws.Range("H122").Value = 2059.1924
$ws.Range("I122").Value = 2064.1667
$ws.Range("K122").Value = 6192.500100000001
$ws.Range("M122").Value = -3742.500100000001

$ws.Range("H132").Value = 1130767.6
$ws.Range("I132").Value = 1249669.5
$ws.Range("K132").Value = 3749008.5
$ws.Range("M132").Value = -3746478.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 442.54544
$ws.Range("I22").Value = 442.54544
$ws.Range("K22").Value = 442.54544
$ws.Range("M22").Value = -269.54544

$ws.Range("H86").Value = 2680
$ws.Range("I86").Value = 2850
$ws.Range("K86").Value = 2850
$ws.Range("M86").Value = -1727

$ws.Range("H89").Value = 2680
$ws.Range("I89").Value = 2850
$ws.Range("K89").Value = 14250
$ws.Range("M89").Value = -8634

$ws.Range("H94").Value = 1879.25
$ws.Range("I94").Value = 1502.3334
$ws.Range("K94").Value = 1502.3334
$ws.Range("M94").Value = -1051.3334

$ws.Range("H134").Value = 2842647.8
$ws.Range("I134").Value = 4766104.5
$ws.Range("K134").Value = 14298313.5
$ws.Range("M134").Value = -14295778.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 765.13336
$ws.Range("I22").Value = 822.53845
$ws.Range("J22").Value = 392
$ws.Range("K22").Value = 822.53845
$ws.Range("L22").Value = 392
$ws.Range("M22").Value = -472.53845
$ws.Range("N22").Value = -1092

$ws.Range("H31").Value = 233470.53
$ws.Range("I31").Value = 391768.12
$ws.Range("K31").Value = 391768.12
$ws.Range("M31").Value = -391473.12

$ws.Range("H34").Value = 233470.53
$ws.Range("I34").Value = 391768.12
$ws.Range("K34").Value = 391768.12
$ws.Range("M34").Value = -391566.12

$ws.Range("H105").Value = 28675.846
$ws.Range("I105").Value = 36318.8
$ws.Range("K105").Value = 36318.8
$ws.Range("M105").Value = -34571.8

$ws.Range("H134").Value = 8816.315000000001
$ws.Range("I134").Value = 10000.5625
$ws.Range("J134").Value = 2500.3333
$ws.Range("K134").Value = 30001.6875
$ws.Range("L134").Value = 7500.999899999999
$ws.Range("M134").Value = -27466.6875
$ws.Range("N134").Value = -12570.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2085.5
$ws.Range("J129").Value = 3999
$ws.Range("L129").Value = 11997
$ws.Range("N129").Value = -21997

$ws.Range("H131").Value = 16018.4
$ws.Range("I131").Value = 678.5714
$ws.Range("J131").Value = 29440.75
$ws.Range("K131").Value = 2035.7142
$ws.Range("L131").Value = 88322.25
$ws.Range("M131").Value = 3004.2858
$ws.Range("N131").Value = -98402.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 4464.4
$ws.Range("I29").Value = 1437.6666
$ws.Range("J29").Value = 9004.5
$ws.Range("K29").Value = 1437.6666
$ws.Range("L29").Value = 9004.5
$ws.Range("M29").Value = -1147.6666
$ws.Range("N29").Value = -9584.5

$ws.Range("H70").Value = 10690.25
$ws.Range("I70").Value = 11525.714
$ws.Range("K70").Value = 11525.714
$ws.Range("M70").Value = -11255.714

$ws.Range("H73").Value = 10690.25
$ws.Range("I73").Value = 11525.714
$ws.Range("K73").Value = 11525.714
$ws.Range("M73").Value = -10589.714

$ws.Range("H122").Value = 55858.4
$ws.Range("I122").Value = 102524.6
$ws.Range("J122").Value = 9192.200000000001
$ws.Range("K122").Value = 307573.8
$ws.Range("L122").Value = 27576.6
$ws.Range("M122").Value = -305123.8
$ws.Range("N122").Value = -32476.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1453.4286
$ws.Range("I22").Value = 665.5714
$ws.Range("J22").Value = 2241.2856
$ws.Range("K22").Value = 665.5714
$ws.Range("L22").Value = 2241.2856
$ws.Range("M22").Value = -370.5714
$ws.Range("N22").Value = -2831.2856

$ws.Range("H27").Value = 1453.4286
$ws.Range("I27").Value = 665.5714
$ws.Range("J27").Value = 2241.2856
$ws.Range("K27").Value = 665.5714
$ws.Range("L27").Value = 2241.2856
$ws.Range("M27").Value = -558.5714
$ws.Range("N27").Value = -2455.2856

$ws.Range("H61").Value = 1598.8182
$ws.Range("I61").Value = 1354.6666
$ws.Range("J61").Value = 2697.5
$ws.Range("K61").Value = 1354.6666
$ws.Range("L61").Value = 2697.5
$ws.Range("M61").Value = -1152.6666
$ws.Range("N61").Value = -3101.5

$ws.Range("H113").Value = 1598.8182
$ws.Range("I113").Value = 1354.6666
$ws.Range("J113").Value = 2697.5
$ws.Range("K113").Value = 1354.6666
$ws.Range("L113").Value = 2697.5
$ws.Range("M113").Value = 815.3334
$ws.Range("N113").Value = -7037.5

$ws.Range("H122").Value = 4857.9585
$ws.Range("I122").Value = 4663.227
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 13989.681
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -11539.681
$ws.Range("N122").Value = -25900

$ws.Range("H136").Value = 51128.652
$ws.Range("I136").Value = 2828.2666
$ws.Range("K136").Value = 8484.799800000001
$ws.Range("M136").Value = -5934.799800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H107").Value = 1700.3438
$ws.Range("I107").Value = 1079.9166
$ws.Range("J107").Value = 3561.625
$ws.Range("K107").Value = 3239.7498
$ws.Range("L107").Value = 10684.875
$ws.Range("M107").Value = -1319.7498
$ws.Range("N107").Value = -14524.875

$ws.Range("H113").Value = 3446.6667
$ws.Range("I113").Value = 1860.625
$ws.Range("J113").Value = 4715.5
$ws.Range("K113").Value = 5581.875
$ws.Range("L113").Value = 14146.5
$ws.Range("M113").Value = -3411.875
$ws.Range("N113").Value = -18486.5
